$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename and update data ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "summ20837374"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 3365.960959972894
$ws.Cells.Item(2, 3).Value = 0.03687366907781617
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 233.4546449101116
$ws.Cells.Item(3, 3).Value = 0.7542595143832623
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 241.6762869141599
$ws.Cells.Item(4, 3).Value = 0.7484621207171274
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -110.0571241271547
$ws.Cells.Item(5, 3).Value = 0.9256075474049701
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -152.1027510165131
$ws.Cells.Item(6, 3).Value = 0.04677873461655273
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -635.525175249887
$ws.Cells.Item(7, 3).Value = 0.001228251691315586
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -35.5770384484283
$ws.Cells.Item(8, 3).Value = [double]"1.108661830855898e-05"
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 237.5786323410834
$ws.Cells.Item(9, 3).Value = 0.0928846606302674
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 1047.903031459733
$ws.Cells.Item(10, 3).Value = [double]"4.585629809145629e-08"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.05874102790562866
$ws.Cells.Item(11, 3).Value = 0.1450446720308035
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"1.320159101608758e-05"
$ws.Cells.Item(12, 3).Value = 0.835798476994565
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 11.06806176936045
$ws.Cells.Item(13, 3).Value = 0.1102042862975298
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 5.30086860485253
$ws.Cells.Item(14, 3).Value = 0.3761278434010181
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -2108.475813750828
$ws.Cells.Item(15, 3).Value = 0.06988320790281026
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1025.961894241101
$ws.Cells.Item(16, 3).Value = 0.2188063261758366
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -28464.2251849918
$ws.Cells.Item(17, 3).Value = 0.0001409283068595271

# --- Sheet 2: rename and update data ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "summ21049838"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5349.128502553171
$ws.Cells.Item(2, 3).Value = 0.001345638506470814
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 277.452132878674
$ws.Cells.Item(3, 3).Value = 0.7156236778794433
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 82.82440668927316
$ws.Cells.Item(4, 3).Value = 0.9143266294478567
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -585.2965304479715
$ws.Cells.Item(5, 3).Value = 0.59245306420361
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -117.0758954999282
$ws.Cells.Item(6, 3).Value = 0.1339312193938089
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -533.5615846909427
$ws.Cells.Item(7, 3).Value = 0.007693096645671414
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -40.29516418430808
$ws.Cells.Item(8, 3).Value = [double]"9.451057504417183e-07"
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 115.3074947441701
$ws.Cells.Item(9, 3).Value = 0.4239306415159082
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 795.1657964600049
$ws.Cells.Item(10, 3).Value = [double]"7.28793775898675e-05"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.1457953299740372
$ws.Cells.Item(11, 3).Value = 0.0004316101405186227
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = 0.0001244806105182545
$ws.Cells.Item(12, 3).Value = 0.06283597645823255
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 0.9717544839608507
$ws.Cells.Item(13, 3).Value = 0.8895286395849195
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 3.897997863482771
$ws.Cells.Item(14, 3).Value = 0.5187344048739018
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -3661.686148783559
$ws.Cells.Item(15, 3).Value = 0.002993382900645289
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -812.121990869991
$ws.Cells.Item(16, 3).Value = 0.3486559415256786
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -29385.25334804377
$ws.Cells.Item(17, 3).Value = 0.0002753476630388929

# --- Sheet 3: rename and update data ---
$ws = $wb.Worksheets.Item(3)
$ws.Name = "summ21299856"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5850.27225990871
$ws.Cells.Item(2, 3).Value = 0.000690615770168713
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 584.6169260414199
$ws.Cells.Item(3, 3).Value = 0.4634712003311845
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 333.5937555496971
$ws.Cells.Item(4, 3).Value = 0.6782657362590563
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -1040.157438042027
$ws.Cells.Item(5, 3).Value = 0.3968677533447218
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -44.84732425123592
$ws.Cells.Item(6, 3).Value = 0.5843808564586476
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -779.021191097972
$ws.Cells.Item(7, 3).Value = 0.0002096513703907842
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -36.36463311778046
$ws.Cells.Item(8, 3).Value = [double]"1.796761479950764e-05"
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 287.910592432921
$ws.Cells.Item(9, 3).Value = 0.05810144536244739
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 1058.080699677677
$ws.Cells.Item(10, 3).Value = [double]"2.57472733408265e-07"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.08210155223067184
$ws.Cells.Item(11, 3).Value = 0.05406990875017941
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"8.818103912941155e-05"
$ws.Cells.Item(12, 3).Value = 0.1820236277583876
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -5.605426983800591
$ws.Cells.Item(13, 3).Value = 0.4406310722145277
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -7.820228556457263
$ws.Cells.Item(14, 3).Value = 0.2229379118811319
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -3549.390606801909
$ws.Cells.Item(15, 3).Value = 0.003920109005800723
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1109.776749655542
$ws.Cells.Item(16, 3).Value = 0.2083285347401958
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -31348.84174140473
$ws.Cells.Item(17, 3).Value = [double]"8.545744693769493e-05"

# --- Sheet 4: rename and update data ---
$ws = $wb.Worksheets.Item(4)
$ws.Name = "summ21561037"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 4498.6178619544
$ws.Cells.Item(2, 3).Value = 0.009211670806643566
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 528.7616585822476
$ws.Cells.Item(3, 3).Value = 0.5057833199143033
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 371.1376123454726
$ws.Cells.Item(4, 3).Value = 0.643233031535251
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -489.6777802069913
$ws.Cells.Item(5, 3).Value = 0.6678065126561344
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -106.4151108394637
$ws.Cells.Item(6, 3).Value = 0.1911270867386392
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -714.1914488245607
$ws.Cells.Item(7, 3).Value = 0.000624900239968115
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -28.8592405306579
$ws.Cells.Item(8, 3).Value = 0.0005544908535299539
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 351.022917397994
$ws.Cells.Item(9, 3).Value = 0.02864917203099888
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 1058.936649705267
$ws.Cells.Item(10, 3).Value = [double]"5.134925448024842e-07"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.09191347949489705
$ws.Cells.Item(11, 3).Value = 0.03099162442879942
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"6.684394453525018e-05"
$ws.Cells.Item(12, 3).Value = 0.3171785692968104
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 0.9509535699179921
$ws.Cells.Item(13, 3).Value = 0.8978924423325962
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -4.146513925919227
$ws.Cells.Item(14, 3).Value = 0.5390067996814498
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -2393.366889116135
$ws.Cells.Item(15, 3).Value = 0.05531488260405852
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1070.720324590863
$ws.Cells.Item(16, 3).Value = 0.2319228835985098
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -26982.10605309781
$ws.Cells.Item(17, 3).Value = 0.0007473893661987739

# --- Sheet 5: rename and update data ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "summ21809169"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5615.922607532159
$ws.Cells.Item(2, 3).Value = 0.0009215764897624678
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 429.6065201387648
$ws.Cells.Item(3, 3).Value = 0.5756736927289766
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 286.8939134736136
$ws.Cells.Item(4, 3).Value = 0.711176201534935
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -626.3609841089382
$ws.Cells.Item(5, 3).Value = 0.5696722431506495
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -117.9917819544622
$ws.Cells.Item(6, 3).Value = 0.1296083541943087
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -676.2951388272393
$ws.Cells.Item(7, 3).Value = 0.000953078726329782
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -30.04720279544197
$ws.Cells.Item(8, 3).Value = 0.0003380044130414919
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 130.7710385677749
$ws.Cells.Item(9, 3).Value = 0.3660656785206662
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 1068.238106466653
$ws.Cells.Item(10, 3).Value = [double]"8.079693014091979e-08"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.1033333040529242
$ws.Cells.Item(11, 3).Value = 0.01183201383989768
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"6.677877527333421e-05"
$ws.Cells.Item(12, 3).Value = 0.3102722712561559
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 1.806812225924554
$ws.Cells.Item(13, 3).Value = 0.7951258239617685
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -4.724487056325254
$ws.Cells.Item(14, 3).Value = 0.4426476712762545
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -4287.887643178815
$ws.Cells.Item(15, 3).Value = 0.0004006730034398722
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1796.653771496048
$ws.Cells.Item(16, 3).Value = 0.03960076194131527
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -25925.69889082369
$ws.Cells.Item(17, 3).Value = 0.001645343715865452

# --- Sheet 6: rename and update data ---
$ws = $wb.Worksheets.Item(6)
$ws.Name = "summ22115737"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5448.524541204495
$ws.Cells.Item(2, 3).Value = 0.001262983784561574
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 593.5773500937255
$ws.Cells.Item(3, 3).Value = 0.439469078115215
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 285.2031468891864
$ws.Cells.Item(4, 3).Value = 0.7126347302218081
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -696.3972075159767
$ws.Cells.Item(5, 3).Value = 0.5270863925381183
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -29.70802880818212
$ws.Cells.Item(6, 3).Value = 0.6986110094738898
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -707.9251650737499
$ws.Cells.Item(7, 3).Value = 0.0004515247841105464
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -29.03798899209144
$ws.Cells.Item(8, 3).Value = 0.000372575045956716
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 185.346927738004
$ws.Cells.Item(9, 3).Value = 0.2035335726929994
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 728.5139078172856
$ws.Cells.Item(10, 3).Value = 0.0003130446310729248
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.06127158538260205
$ws.Cells.Item(11, 3).Value = 0.1353845468165388
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"2.708660597720189e-05"
$ws.Cells.Item(12, 3).Value = 0.6788725962741775
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -0.6733593874746493
$ws.Cells.Item(13, 3).Value = 0.92468183381956
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 1.516715076588317
$ws.Cells.Item(14, 3).Value = 0.8091082024188906
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -3760.907412413209
$ws.Cells.Item(15, 3).Value = 0.001751816246989149
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1817.554109340256
$ws.Cells.Item(16, 3).Value = 0.03766461029595138
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -26085.72285950735
$ws.Cells.Item(17, 3).Value = 0.0005613644395418629

# --- Sheet 7: rename and update data ---
$ws = $wb.Worksheets.Item(7)
$ws.Name = "summ22413402"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 5944.299767717822
$ws.Cells.Item(2, 3).Value = 0.001092051258705421
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 212.6180055041206
$ws.Cells.Item(3, 3).Value = 0.8171116828646299
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 240.0687262637486
$ws.Cells.Item(4, 3).Value = 0.7956071611923895
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -1119.713664001291
$ws.Cells.Item(5, 3).Value = 0.3959149970217603
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -72.85178579031711
$ws.Cells.Item(6, 3).Value = 0.3822844342267643
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -486.0801859712703
$ws.Cells.Item(7, 3).Value = 0.02157027187399206
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -37.97141597774628
$ws.Cells.Item(8, 3).Value = [double]"1.590029829348337e-05"
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 56.76294460992321
$ws.Cells.Item(9, 3).Value = 0.7085884437027643
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 937.1351050018233
$ws.Cells.Item(10, 3).Value = [double]"5.702839741853358e-06"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.1093580851006513
$ws.Cells.Item(11, 3).Value = 0.01021036341452951
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"9.930527155781886e-05"
$ws.Cells.Item(12, 3).Value = 0.1316202870954075
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -1.260011352783014
$ws.Cells.Item(13, 3).Value = 0.8638171655669883
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -0.2716553571941285
$ws.Cells.Item(14, 3).Value = 0.965995750943422
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -4422.175895304332
$ws.Cells.Item(15, 3).Value = 0.0003858584441325532
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -1692.335241506312
$ws.Cells.Item(16, 3).Value = 0.06015653672952342
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -30673.26673135663
$ws.Cells.Item(17, 3).Value = 0.0001081941007922935

# --- Sheet 8: rename and update data ---
$ws = $wb.Worksheets.Item(8)
$ws.Name = "summ22681308"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 4906.998284235458
$ws.Cells.Item(2, 3).Value = 0.009221924183737873
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = 907.8392396287479
$ws.Cells.Item(3, 3).Value = 0.3496967866200875
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = 743.4014703546682
$ws.Cells.Item(4, 3).Value = 0.44731618894263
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -193.2166165679532
$ws.Cells.Item(5, 3).Value = 0.8888647162049761
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -24.01642725074973
$ws.Cells.Item(6, 3).Value = 0.774872751926734
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -682.7617756419327
$ws.Cells.Item(7, 3).Value = 0.001426367189460273
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -31.78325856016352
$ws.Cells.Item(8, 3).Value = 0.0002613117686733499
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 140.7796217444195
$ws.Cells.Item(9, 3).Value = 0.3648298027839596
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 954.2834535453119
$ws.Cells.Item(10, 3).Value = [double]"7.306133696148864e-06"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.103266876978487
$ws.Cells.Item(11, 3).Value = 0.0162245471822705
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"5.80863175453272e-05"
$ws.Cells.Item(12, 3).Value = 0.3975974299659188
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = 4.404273513102693
$ws.Cells.Item(13, 3).Value = 0.5551956017394712
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = 0.1042461641649268
$ws.Cells.Item(14, 3).Value = 0.9884426196342359
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -4081.966986649226
$ws.Cells.Item(15, 3).Value = 0.001317930678391426
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -2241.683958949222
$ws.Cells.Item(16, 3).Value = 0.01728719990973466
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -30317.9811760569
$ws.Cells.Item(17, 3).Value = 0.0003128674532824585

# --- Sheet 9: rename and update data ---
$ws = $wb.Worksheets.Item(9)
$ws.Name = "summ22949014"
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(2, 1).Value = "Intercept"
$ws.Cells.Item(2, 2).Value = 4236.368676822472
$ws.Cells.Item(2, 3).Value = 0.02443034508435456
$ws.Cells.Item(3, 1).Value = "Education[T.Secondary]"
$ws.Cells.Item(3, 2).Value = -161.2890478174619
$ws.Cells.Item(3, 3).Value = 0.8744968169592915
$ws.Cells.Item(4, 1).Value = "Education[T.University]"
$ws.Cells.Item(4, 2).Value = -83.26868465970216
$ws.Cells.Item(4, 3).Value = 0.9354533978126276
$ws.Cells.Item(5, 1).Value = "Education[T.Unknown/Other]"
$ws.Cells.Item(5, 2).Value = -1046.999687885616
$ws.Cells.Item(5, 3).Value = 0.4171133480026011
$ws.Cells.Item(6, 1).Value = "HHSize"
$ws.Cells.Item(6, 2).Value = -71.39354240439202
$ws.Cells.Item(6, 3).Value = 0.3682450529399883
$ws.Cells.Item(7, 1).Value = "Sex"
$ws.Cells.Item(7, 2).Value = -793.4371509376535
$ws.Cells.Item(7, 3).Value = 0.0001027776964294664
$ws.Cells.Item(8, 1).Value = "Age"
$ws.Cells.Item(8, 2).Value = -26.94521777416928
$ws.Cells.Item(8, 3).Value = 0.001113969759926302
$ws.Cells.Item(9, 1).Value = "DistSubcenter_res"
$ws.Cells.Item(9, 2).Value = 422.2693396540243
$ws.Cells.Item(9, 3).Value = 0.004046697439039407
$ws.Cells.Item(10, 1).Value = "DistCenter_res"
$ws.Cells.Item(10, 2).Value = 1048.070017847687
$ws.Cells.Item(10, 3).Value = [double]"1.878808364234321e-07"
$ws.Cells.Item(11, 1).Value = "UrbPopDensity_res"
$ws.Cells.Item(11, 2).Value = -0.07729967600167534
$ws.Cells.Item(11, 3).Value = 0.06167735509582199
$ws.Cells.Item(12, 1).Value = "UrbBuildDensity_res"
$ws.Cells.Item(12, 2).Value = [double]"6.739373414248683e-05"
$ws.Cells.Item(12, 3).Value = 0.3122107454650849
$ws.Cells.Item(13, 1).Value = "IntersecDensity_res"
$ws.Cells.Item(13, 2).Value = -1.230517088802156
$ws.Cells.Item(13, 3).Value = 0.8658418697415887
$ws.Cells.Item(14, 1).Value = "street_length_res"
$ws.Cells.Item(14, 2).Value = -1.507272330718672
$ws.Cells.Item(14, 3).Value = 0.8264184487893071
$ws.Cells.Item(15, 1).Value = "LU_Comm_res"
$ws.Cells.Item(15, 2).Value = -2696.07855821297
$ws.Cells.Item(15, 3).Value = 0.02545531837714374
$ws.Cells.Item(16, 1).Value = "LU_UrbFab_res"
$ws.Cells.Item(16, 2).Value = -247.4090621087116
$ws.Cells.Item(16, 3).Value = 0.7804184593798239
$ws.Cells.Item(17, 1).Value = "bike_lane_share_res"
$ws.Cells.Item(17, 2).Value = -24232.11084501407
$ws.Cells.Item(17, 3).Value = 0.00148759920555563
